$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New gerber export shifted a couple of existing placements too:
# Row 4 (PROG header) gets new coordinates
$ws.Range("D4").Value = "59.69mm"
$ws.Range("E4").Value = "40.767mm"
$ws.Range("F4").Value = "59.69mm"
$ws.Range("G4").Value = "40.767mm"
$ws.Range("H4").Value = "59.69mm"
$ws.Range("I4").Value = "34.417mm"

# Row 6: the old separate "5V" and "12V" header rows merge into one "12V" row here
$ws.Range("A6").Value = "12V"
$ws.Range("D6").Value = "52.832mm"
$ws.Range("E6").Value = "22.098mm"
$ws.Range("F6").Value = "52.832mm"
$ws.Range("G6").Value = "22.098mm"
$ws.Range("H6").Value = "52.832mm"
$ws.Range("I6").Value = "22.098mm"
$ws.Range("N6").Value = "PZ254V-11-01P"

# Row 7 -> C1 (100nF capacitor), was the "PROG" header row content shifted up in the diff
$ws.Range("A7").Value = "C1"
$ws.Range("B7").Value = "MES104J2A-7-50R0"
$ws.Range("C7").Value = "CAP-TH_L7.2-W4.0-P5.00-D0.5"
$ws.Range("D7").Value = "43.561mm"
$ws.Range("E7").Value = "27.305mm"
$ws.Range("F7").Value = "43.561mm"
$ws.Range("G7").Value = "27.305mm"
$ws.Range("H7").Value = "41.062mm"
$ws.Range("I7").Value = "27.305mm"
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = "T"
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = "No"
$ws.Range("N7").Value = "100nF"

# Row 8 -> C2
$ws.Range("A8").Value = "C2"
$ws.Range("B8").Value = "MES104J2A-7-50R0"
$ws.Range("C8").Value = "CAP-TH_L7.2-W4.0-P5.00-D0.5"
$ws.Range("D8").Value = "53.467mm"
$ws.Range("E8").Value = "27.305mm"
$ws.Range("F8").Value = "53.467mm"
$ws.Range("G8").Value = "27.305mm"
$ws.Range("H8").Value = "55.966mm"
$ws.Range("I8").Value = "27.305mm"
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = "T"
$ws.Range("L8").Value = 180
$ws.Range("M8").Value = "No"
$ws.Range("N8").Value = "100nF"

# Row 9 -> U7
$ws.Range("A9").Value = "U7"
$ws.Range("B9").Value = "7805_C305416"
$ws.Range("C9").Value = "TO-220-3_L10.4-W15.5-P3.00"
$ws.Range("D9").Value = "49.149mm"
$ws.Range("E9").Value = "31.877mm"
$ws.Range("F9").Value = "49.149mm"
$ws.Range("G9").Value = "31.877mm"
$ws.Range("H9").Value = "46.149mm"
$ws.Range("I9").Value = "31.877mm"
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = "T"
$ws.Range("L9").Value = 270
$ws.Range("M9").Value = "No"
# "7805" looks numeric, so force it to be stored as text (matches the source PnP export)
$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "7805"
$ws.Range("N9").ClearFormats()

# Row 10 -> U3 (Level Shifter), replaces old row 10 content
$ws.Range("A10").Value = "U3"
$ws.Range("B10").Value = "Level Shifter"
$ws.Range("C10").Value = "Level Shifter Footprint"
$ws.Range("D10").Value = "32.15mm"
$ws.Range("E10").Value = "35.876mm"
$ws.Range("F10").Value = "25.781mm"
$ws.Range("G10").Value = "43.942mm"
$ws.Range("H10").Value = "27.051mm"
$ws.Range("I10").Value = "42.226mm"
$ws.Range("J10").Value = 12
$ws.Range("K10").Value = "T"
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = "No"
$ws.Range("N10").Value = "Level Shifter"

# The table shrank by one row (old row 11 is no longer needed) - remove it
$ws.Rows.Item(11).Delete()
